$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "foo"
$ws.Range("C2").Value = "bar"
$ws.Range("D2").Value = "ok"
$ws.Range("E2").Value = "haha"
$ws.Range("C3").Value = "foo foo"
$ws.Range("F3").Value = "bar bar"
